# Natmi following Dr Hou advice
# Rebuild the LR-pair result table: each sending cluster (ECs, FAPs, M1, M2,
# Neutro, sCs) now gets TWO rows (Target cluster = FAPs, then sCs) instead of
# a single row (Target cluster = sCs only), with recomputed statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Lpar3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 106.8785896666667
$ws.Range("H2").Value = 320.635769
$ws.Range("I2").Value = 0.1508748302900445
$ws.Range("J2").Value = 0.1508748302900445
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09535966666666668
$ws.Range("N2").Value = 0.286079
$ws.Range("O2").Value = 0.08150805239577472
$ws.Range("P2").Value = 0.08150805239577472
$ws.Range("Q2").Value = 10.19190668441678
$ws.Range("R2").Value = 91.72716015975101
$ws.Range("S2").Value = 0.01229751357248457
$ws.Range("T2").Value = 0.01229751357248457

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Lpar3"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 106.8785896666667
$ws.Range("H3").Value = 320.635769
$ws.Range("I3").Value = 0.1508748302900445
$ws.Range("J3").Value = 0.1508748302900445
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.074582
$ws.Range("N3").Value = 3.223746
$ws.Range("O3").Value = 0.9184919476042253
$ws.Range("P3").Value = 0.9184919476042253
$ws.Range("Q3").Value = 114.849808641186
$ws.Range("R3").Value = 1033.648277770674
$ws.Range("S3").Value = 0.13857731671756
$ws.Range("T3").Value = 0.13857731671756

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Lpar3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 77.232923
$ws.Range("H4").Value = 231.698769
$ws.Range("I4").Value = 0.1090256166999485
$ws.Range("J4").Value = 0.1090256166999485
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.09535966666666668
$ws.Range("N4").Value = 0.286079
$ws.Range("O4").Value = 0.08150805239577472
$ws.Range("P4").Value = 0.08150805239577472
$ws.Range("Q4").Value = 7.364905792972334
$ws.Range("R4").Value = 66.284152136751
$ws.Range("S4").Value = 0.008886465678461054
$ws.Range("T4").Value = 0.008886465678461054

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Lpar3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 77.232923
$ws.Range("H5").Value = 231.698769
$ws.Range("I5").Value = 0.1090256166999485
$ws.Range("J5").Value = 0.1090256166999485
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.074582
$ws.Range("N5").Value = 3.223746
$ws.Range("O5").Value = 0.9184919476042253
$ws.Range("P5").Value = 0.9184919476042253
$ws.Range("Q5").Value = 82.99310886318601
$ws.Range("R5").Value = 746.937979768674
$ws.Range("S5").Value = 0.1001391510214874
$ws.Range("T5").Value = 0.1001391510214875

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Lpar3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 136.676337
$ws.Range("H6").Value = 410.029011
$ws.Range("I6").Value = 0.1929387280825172
$ws.Range("J6").Value = 0.1929387280825172
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.09535966666666668
$ws.Range("N6").Value = 0.286079
$ws.Range("O6").Value = 0.08150805239577472
$ws.Range("P6").Value = 0.08150805239577472
$ws.Range("Q6").Value = 13.033409937541
$ws.Range("R6").Value = 117.300689437869
$ws.Range("S6").Value = 0.01572605995772395
$ws.Range("T6").Value = 0.01572605995772395

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Lpar3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 136.676337
$ws.Range("H7").Value = 410.029011
$ws.Range("I7").Value = 0.1929387280825172
$ws.Range("J7").Value = 0.1929387280825172
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.074582
$ws.Range("N7").Value = 3.223746
$ws.Range("O7").Value = 0.9184919476042253
$ws.Range("P7").Value = 0.9184919476042253
$ws.Range("Q7").Value = 146.869931566134
$ws.Range("R7").Value = 1321.829384095206
$ws.Range("S7").Value = 0.1772126681247933
$ws.Range("T7").Value = 0.1772126681247933

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Lpar3"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 143.4557186666667
$ws.Range("H8").Value = 430.367156
$ws.Range("I8").Value = 0.2025088212285795
$ws.Range("J8").Value = 0.2025088212285795
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.09535966666666668
$ws.Range("N8").Value = 0.286079
$ws.Range("O8").Value = 0.08150805239577472
$ws.Range("P8").Value = 0.08150805239577472
$ws.Range("Q8").Value = 13.67988951348045
$ws.Range("R8").Value = 123.119005621324
$ws.Range("S8").Value = 0.01650609961130563
$ws.Range("T8").Value = 0.01650609961130563

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Lpar3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 143.4557186666667
$ws.Range("H9").Value = 430.367156
$ws.Range("I9").Value = 0.2025088212285795
$ws.Range("J9").Value = 0.2025088212285795
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.074582
$ws.Range("N9").Value = 3.223746
$ws.Range("O9").Value = 0.9184919476042253
$ws.Range("P9").Value = 0.9184919476042253
$ws.Range("Q9").Value = 154.154933076264
$ws.Range("R9").Value = 1387.394397686376
$ws.Range("S9").Value = 0.1860027216172739
$ws.Range("T9").Value = 0.1860027216172739

# Row 10
$ws.Range("A10").Value = "Neutro"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "Lpar3"
$ws.Range("D10").Value = "FAPs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 198.5982106666667
$ws.Range("H10").Value = 595.794632
$ws.Range("I10").Value = 0.2803505493821544
$ws.Range("J10").Value = 0.2803505493821544
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.09535966666666668
$ws.Range("N10").Value = 0.286079
$ws.Range("O10").Value = 0.08150805239577472
$ws.Range("P10").Value = 0.08150805239577472
$ws.Range("Q10").Value = 18.93825916976978
$ws.Range("R10").Value = 170.444332527928
$ws.Range("S10").Value = 0.02285082726822487
$ws.Range("T10").Value = 0.02285082726822487

# Row 11
$ws.Range("A11").Value = "Neutro"
$ws.Range("B11").Value = "Gnai2"
$ws.Range("C11").Value = "Lpar3"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 198.5982106666667
$ws.Range("H11").Value = 595.794632
$ws.Range("I11").Value = 0.2803505493821544
$ws.Range("J11").Value = 0.2803505493821544
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.074582
$ws.Range("N11").Value = 3.223746
$ws.Range("O11").Value = 0.9184919476042253
$ws.Range("P11").Value = 0.9184919476042253
$ws.Range("Q11").Value = 213.410062414608
$ws.Range("R11").Value = 1920.690561731472
$ws.Range("S11").Value = 0.2574997221139295
$ws.Range("T11").Value = 0.2574997221139295

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Gnai2"
$ws.Range("C12").Value = "Lpar3"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 45.55066433333334
$ws.Range("H12").Value = 136.651993
$ws.Range("I12").Value = 0.06430145431675577
$ws.Range("J12").Value = 0.06430145431675577
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.09535966666666668
$ws.Range("N12").Value = 0.286079
$ws.Range("O12").Value = 0.08150805239577472
$ws.Range("P12").Value = 0.08150805239577472
$ws.Range("Q12").Value = 4.34369616727189
$ws.Range("R12").Value = 39.09326550544701
$ws.Range("S12").Value = 0.005241086307574644
$ws.Range("T12").Value = 0.005241086307574644

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Gnai2"
$ws.Range("C13").Value = "Lpar3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 45.55066433333334
$ws.Range("H13").Value = 136.651993
$ws.Range("I13").Value = 0.06430145431675577
$ws.Range("J13").Value = 0.06430145431675577
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.074582
$ws.Range("N13").Value = 3.223746
$ws.Range("O13").Value = 0.9184919476042253
$ws.Range("P13").Value = 0.9184919476042253
$ws.Range("Q13").Value = 48.94792398064201
$ws.Range("R13").Value = 440.5313158257781
$ws.Range("S13").Value = 0.05906036800918112
$ws.Range("T13").Value = 0.05906036800918112
